$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 3, pushing existing rows 3-7 down to 4-8.
$ws.Rows.Item(3).Insert()

# Populate the new row 3 with values (copy the constant columns, then set the changed ones).
$ws.Range("A3").Value = 11
$ws.Range("B3").Value = "Vega Monumental Concepción"
$ws.Range("C3").Value = "Bíobío"
$ws.Range("D3").Value = 44742
$ws.Range("E3").Value = 8
$ws.Range("F3").Value = "Fruta"
$ws.Range("G3").Value = 100107
$ws.Range("H3").Value = "Otros"
$ws.Range("I3").Value = 100107001
$ws.Range("J3").Value = "Caqui"
$ws.Range("K3").Value = "Mankaki"
$ws.Range("L3").Value = "Segunda"
$ws.Range("M3").Value = 100
$ws.Range("N3").Value = 14000
$ws.Range("O3").Value = 15000
$ws.Range("P3").Value = 14500
$ws.Range("Q3").Value = "$/caja 18 kilos granel"
$ws.Range("R3").Value = "Región de O'Higgins"
$ws.Range("S3").Value = 806
$ws.Range("T3").Value = 18

# Match the date format style used in column D (s="2")
$ws.Range("D3").NumberFormat = "YYYY-MM-DD HH:MM:SS"
